$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.80646166666667
$ws.Range("H2").Value = 41.419385
$ws.Range("I2").Value = 0.2210624443376167
$ws.Range("J2").Value = 0.2210624443376167
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.525001
$ws.Range("N2").Value = 1.575003
$ws.Range("O2").Value = 0.07505143515225263
$ws.Range("P2").Value = 0.07505143515225261
$ws.Range("Q2").Value = 7.248406181461667
$ws.Range("R2").Value = 65.235655633155
$ws.Range("S2").Value = 0.0165910537058031
$ws.Range("T2").Value = 0.01659105370580309

$ws.Range("G3").Value = 13.80646166666667
$ws.Range("H3").Value = 41.419385
$ws.Range("I3").Value = 0.2210624443376167
$ws.Range("J3").Value = 0.2210624443376167
$ws.Range("M3").Value = 4.457871000000001
$ws.Range("O3").Value = 0.6372742457130702
$ws.Range("P3").Value = 0.6372742457130701
$ws.Range("Q3").Value = 61.54742507644501
$ws.Range("R3").Value = 553.9268256880051
$ws.Range("S3").Value = 0.1408774024707423
$ws.Range("T3").Value = 0.1408774024707423

$ws.Range("G4").Value = 13.80646166666667
$ws.Range("H4").Value = 41.419385
$ws.Range("I4").Value = 0.2210624443376167
$ws.Range("J4").Value = 0.2210624443376167
$ws.Range("M4").Value = 2.012344
$ws.Range("N4").Value = 6.037032
$ws.Range("O4").Value = 0.2876743191346772
$ws.Range("P4").Value = 0.2876743191346771
$ws.Range("Q4").Value = 27.78335029614666
$ws.Range("R4").Value = 250.05015266532
$ws.Range("S4").Value = 0.06359398816107137
$ws.Range("T4").Value = 0.06359398816107135

$ws.Range("I5").Value = 0.1985406876033273
$ws.Range("J5").Value = 0.1985406876033273
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.525001
$ws.Range("N5").Value = 1.575003
$ws.Range("O5").Value = 0.07505143515225263
$ws.Range("P5").Value = 0.07505143515225261
$ws.Range("Q5").Value = 6.509941349864667
$ws.Range("R5").Value = 58.589472148782
$ws.Range("S5").Value = 0.01490076354074477
$ws.Range("T5").Value = 0.01490076354074477

$ws.Range("I6").Value = 0.1985406876033273
$ws.Range("J6").Value = 0.1985406876033273
$ws.Range("M6").Value = 4.457871000000001
$ws.Range("O6").Value = 0.6372742457130702
$ws.Range("P6").Value = 0.6372742457130701
$ws.Range("Q6").Value = 55.27699710145801
$ws.Range("R6").Value = 497.4929739131221
$ws.Range("S6").Value = 0.1265248669357647
$ws.Range("T6").Value = 0.1265248669357647

$ws.Range("I7").Value = 0.1985406876033273
$ws.Range("J7").Value = 0.1985406876033273
$ws.Range("M7").Value = 2.012344
$ws.Range("N7").Value = 6.037032
$ws.Range("O7").Value = 0.2876743191346772
$ws.Range("P7").Value = 0.2876743191346771
$ws.Range("Q7").Value = 24.95279326277867
$ws.Range("R7").Value = 224.575139365008
$ws.Range("S7").Value = 0.05711505712681783
$ws.Range("T7").Value = 0.05711505712681782

$ws.Range("G8").Value = 20.11717366666667
$ws.Range("H8").Value = 60.35152100000001
$ws.Range("I8").Value = 0.3221065390457394
$ws.Range("J8").Value = 0.3221065390457393
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.525001
$ws.Range("N8").Value = 1.575003
$ws.Range("O8").Value = 0.07505143515225263
$ws.Range("P8").Value = 0.07505143515225261
$ws.Range("Q8").Value = 10.56153629217367
$ws.Range("R8").Value = 95.05382662956302
$ws.Range("S8").Value = 0.02417455802730784
$ws.Range("T8").Value = 0.02417455802730783

$ws.Range("G9").Value = 20.11717366666667
$ws.Range("H9").Value = 60.35152100000001
$ws.Range("I9").Value = 0.3221065390457394
$ws.Range("J9").Value = 0.3221065390457393
$ws.Range("M9").Value = 4.457871000000001
$ws.Range("O9").Value = 0.6372742457130702
$ws.Range("P9").Value = 0.6372742457130701
$ws.Range("Q9").Value = 89.67976509059703
$ws.Range("R9").Value = 807.1178858153733
$ws.Range("S9").Value = 0.2052702017096212
$ws.Range("T9").Value = 0.2052702017096211

$ws.Range("G10").Value = 20.11717366666667
$ws.Range("H10").Value = 60.35152100000001
$ws.Range("I10").Value = 0.3221065390457394
$ws.Range("J10").Value = 0.3221065390457393
$ws.Range("M10").Value = 2.012344
$ws.Range("N10").Value = 6.037032
$ws.Range("O10").Value = 0.2876743191346772
$ws.Range("P10").Value = 0.2876743191346771
$ws.Range("Q10").Value = 40.48267372507468
$ws.Range("R10").Value = 364.344063525672
$ws.Range("S10").Value = 0.09266177930881038
$ws.Range("T10").Value = 0.09266177930881035

$ws.Range("G11").Value = 4.066281666666666
$ws.Range("H11").Value = 12.198845
$ws.Range("I11").Value = 0.06510735236159866
$ws.Range("J11").Value = 0.06510735236159866
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.525001
$ws.Range("N11").Value = 1.575003
$ws.Range("O11").Value = 0.07505143515225263
$ws.Range("P11").Value = 0.07505143515225261
$ws.Range("Q11").Value = 2.134801941281666
$ws.Range("R11").Value = 19.213217471535
$ws.Range("S11").Value = 0.004886400233701384
$ws.Range("T11").Value = 0.004886400233701383

$ws.Range("G12").Value = 4.066281666666666
$ws.Range("H12").Value = 12.198845
$ws.Range("I12").Value = 0.06510735236159866
$ws.Range("J12").Value = 0.06510735236159866
$ws.Range("M12").Value = 4.457871000000001
$ws.Range("O12").Value = 0.6372742457130702
$ws.Range("P12").Value = 0.6372742457130701
$ws.Range("Q12").Value = 18.126959119665
$ws.Range("R12").Value = 163.142632076985
$ws.Range("S12").Value = 0.04149123886661287
$ws.Range("T12").Value = 0.04149123886661286

$ws.Range("G13").Value = 4.066281666666666
$ws.Range("H13").Value = 12.198845
$ws.Range("I13").Value = 0.06510735236159866
$ws.Range("J13").Value = 0.06510735236159866
$ws.Range("M13").Value = 2.012344
$ws.Range("N13").Value = 6.037032
$ws.Range("O13").Value = 0.2876743191346772
$ws.Range("P13").Value = 0.2876743191346771
$ws.Range("Q13").Value = 8.182757514226665
$ws.Range("R13").Value = 73.64481762803999
$ws.Range("S13").Value = 0.01872971326128441
$ws.Range("T13").Value = 0.01872971326128441

$ws.Range("G14").Value = 12.06524866666667
$ws.Range("H14").Value = 36.195746
$ws.Range("I14").Value = 0.1931829766517179
$ws.Range("J14").Value = 0.1931829766517179
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.525001
$ws.Range("N14").Value = 1.575003
$ws.Range("O14").Value = 0.07505143515225263
$ws.Range("P14").Value = 0.07505143515225261
$ws.Range("Q14").Value = 6.334267615248668
$ws.Range("R14").Value = 57.008408537238
$ws.Range("S14").Value = 0.01449865964469554
$ws.Range("T14").Value = 0.01449865964469554

$ws.Range("G15").Value = 12.06524866666667
$ws.Range("H15").Value = 36.195746
$ws.Range("I15").Value = 0.1931829766517179
$ws.Range("J15").Value = 0.1931829766517179
$ws.Range("M15").Value = 4.457871000000001
$ws.Range("O15").Value = 0.6372742457130702
$ws.Range("P15").Value = 0.6372742457130701
$ws.Range("Q15").Value = 53.78532213892201
$ws.Range("R15").Value = 484.0678992502981
$ws.Range("S15").Value = 0.1231105357303292
$ws.Range("T15").Value = 0.1231105357303292

$ws.Range("G16").Value = 12.06524866666667
$ws.Range("H16").Value = 36.195746
$ws.Range("I16").Value = 0.1931829766517179
$ws.Range("J16").Value = 0.1931829766517179
$ws.Range("M16").Value = 2.012344
$ws.Range("N16").Value = 6.037032
$ws.Range("O16").Value = 0.2876743191346772
$ws.Range("P16").Value = 0.2876743191346771
$ws.Range("Q16").Value = 24.27943076287467
$ws.Range("R16").Value = 218.514876865872
$ws.Range("S16").Value = 0.05557378127669318
$ws.Range("T16").Value = 0.05557378127669316
